# AShot without inspirepak and jenkins propfile
#
# Re-labels the "Ink / Varnish" / "Plate" / "Sheet" material rows on the
# "Material" sheet so that the Element column (B) - and the dependent
# Process/Measurement columns that travel with a given "Folded Sheet" /
# "Folded Sheet 1" / "Folded Sheet 2" group - line up the way the refreshed
# job-material export expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Folded Sheet  4p" block (was "Folded Sheet 2  4p") ---------------
$ws.Cells.Item(2,2).Value = "Folded Sheet  4p"
$ws.Cells.Item(2,4).Value = "Black - Sheet-fed Offset - "
$ws.Cells.Item(2,5).Value = "'0.07"
$ws.Cells.Item(2,8).Value = "M594GK - Black (General)"

$ws.Cells.Item(3,2).Value = "Folded Sheet  4p"
$ws.Cells.Item(3,4).Value = "Yellow - Sheet-fed Offset - "
$ws.Cells.Item(3,5).Value = "'0.07"

$ws.Cells.Item(4,2).Value = "Folded Sheet  4p"
$ws.Cells.Item(4,5).Value = "'0.07"

$ws.Cells.Item(5,2).Value = "Folded Sheet  4p"
$ws.Cells.Item(5,4).Value = "Magenta - Sheet-fed Offset - "
$ws.Cells.Item(5,5).Value = "'0.07"
$ws.Cells.Item(5,8).Value = "M594GN - 4/C Process (General)"

# --- "Folded Sheet 1  4p" block (was "Folded Sheet  4p") ---------------
$ws.Cells.Item(6,2).Value = "Folded Sheet 1  4p"
$ws.Cells.Item(6,4).Value = "Yellow - Sheet-fed Offset - "

$ws.Cells.Item(7,2).Value = "Folded Sheet 1  4p"

$ws.Cells.Item(8,2).Value = "Folded Sheet 1  4p"
$ws.Cells.Item(8,4).Value = "Magenta - Sheet-fed Offset - "

$ws.Cells.Item(9,2).Value = "Folded Sheet 1  4p"
$ws.Cells.Item(9,4).Value = "Cyan - Sheet-fed Offset - "

# --- "Folded Sheet 2  4p" block (was "Folded Sheet 1  4p") -------------
$ws.Cells.Item(10,2).Value = "Folded Sheet 2  4p"
$ws.Cells.Item(10,4).Value = "Yellow - Sheet-fed Offset - "
$ws.Cells.Item(10,5).Value = "'0.08"

$ws.Cells.Item(11,2).Value = "Folded Sheet 2  4p"
$ws.Cells.Item(11,5).Value = "'0.08"

$ws.Cells.Item(12,2).Value = "Folded Sheet 2  4p"
$ws.Cells.Item(12,4).Value = "Magenta - Sheet-fed Offset - "
$ws.Cells.Item(12,5).Value = "'0.08"

$ws.Cells.Item(13,2).Value = "Folded Sheet 2  4p"
$ws.Cells.Item(13,5).Value = "'0.08"

# --- Plate rows: re-point each plate to its matching folded sheet ------
$ws.Cells.Item(14,2).Value = "Plate - Folded Sheet 2  4p"
$ws.Cells.Item(15,2).Value = "Plate - Folded Sheet  4p"
$ws.Cells.Item(16,2).Value = "Plate - Folded Sheet 1  4p"

# --- Sheet (paper) rows: re-point + correct the Per-M quantities -------
$ws.Cells.Item(17,2).Value = "Folded Sheet  4p"
$ws.Cells.Item(17,5).Value = "'1,001.00"

$ws.Cells.Item(18,2).Value = "Folded Sheet 1  4p"
$ws.Cells.Item(18,5).Value = "'1,106.00"

$ws.Cells.Item(19,2).Value = "Folded Sheet 2  4p"
$ws.Cells.Item(19,5).Value = "'1,213.00"
